$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 81.58
$ws.Range("I15").Value = 81.58
$ws.Range("K15").Value = 244.74
$ws.Range("M15").Value = -75.74000000000001
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H121").Value = 2225.5417
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 2225.5417
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 6676.625100000001
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -10170.6251
$ws.Range("H129").Value = 941.63635
$ws.Range("I129").Value = 448.5
$ws.Range("J129").Value = 957.0469000000001
$ws.Range("K129").Value = 1345.5
$ws.Range("L129").Value = 2871.1407
$ws.Range("M129").Value = 3654.5
$ws.Range("N129").Value = -12871.1407
$ws.Range("H133").Value = 51384.832
$ws.Range("J133").Value = 51384.832
$ws.Range("L133").Value = 51384.832
$ws.Range("N133").Value = -61504.832
$ws.Range("H137").Value = 5612.5
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 2696.7356
$ws.Range("I138").Value = 1822
$ws.Range("J138").Value = 2909.1714
$ws.Range("K138").Value = 5466
$ws.Range("L138").Value = 8727.514200000001
$ws.Range("M138").Value = -326
$ws.Range("N138").Value = -19007.5142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11414.75
$ws.Range("I32").Value = 7472.1396
$ws.Range("K32").Value = 7472.1396
$ws.Range("M32").Value = -7185.1396
$ws.Range("H63").Value = 9896421
$ws.Range("I63").Value = 11545025
$ws.Range("J63").Value = 4800
$ws.Range("K63").Value = 11545025
$ws.Range("L63").Value = 4800
$ws.Range("M63").Value = -11544339
$ws.Range("N63").Value = -6172
$ws.Range("H66").Value = 9896421
$ws.Range("I66").Value = 11545025
$ws.Range("J66").Value = 4800
$ws.Range("K66").Value = 57725125
$ws.Range("L66").Value = 24000
$ws.Range("M66").Value = -57721693
$ws.Range("N66").Value = -30864
$ws.Range("H74").Value = 1696.8948
$ws.Range("I74").Value = 1197.0605
$ws.Range("J74").Value = 4995.8
$ws.Range("K74").Value = 1197.0605
$ws.Range("L74").Value = 4995.8
$ws.Range("M74").Value = -323.0605
$ws.Range("N74").Value = -6743.8
$ws.Range("H77").Value = 1696.8948
$ws.Range("I77").Value = 1197.0605
$ws.Range("J77").Value = 4995.8
$ws.Range("K77").Value = 5985.3025
$ws.Range("L77").Value = 24979
$ws.Range("M77").Value = -1617.3025
$ws.Range("N77").Value = -33715
$ws.Range("H98").Value = 59800
$ws.Range("J98").Value = 59800
$ws.Range("L98").Value = 59800
$ws.Range("N98").Value = -65790
$ws.Range("H103").Value = 34285.715
$ws.Range("J103").Value = 34285.715
$ws.Range("L103").Value = 34285.715
$ws.Range("N103").Value = -36629.715
$ws.Range("H122").Value = 3299.1667
$ws.Range("I122").Value = 1768.6666
$ws.Range("J122").Value = 4829.6665
$ws.Range("K122").Value = 5305.9998
$ws.Range("L122").Value = 14488.9995
$ws.Range("M122").Value = -2855.9998
$ws.Range("N122").Value = -19388.9995
$ws.Range("H132").Value = 2397.0232
$ws.Range("I132").Value = 1151.3182
$ws.Range("J132").Value = 3702.0476
$ws.Range("K132").Value = 3453.9546
$ws.Range("L132").Value = 11106.1428
$ws.Range("M132").Value = -923.9546
$ws.Range("N132").Value = -16166.1428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3092.1738
$ws.Range("I134").Value = 1782.2354
$ws.Range("J134").Value = 6803.6665
$ws.Range("K134").Value = 5346.706200000001
$ws.Range("L134").Value = 20410.9995
$ws.Range("M134").Value = -2811.706200000001
$ws.Range("N134").Value = -25480.9995
$ws.Range("H137").Value = 33109.285
$ws.Range("J137").Value = 33109.285
$ws.Range("L137").Value = 33109.285
$ws.Range("N137").Value = -43309.285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4832811
$ws.Range("I16").Value = 10101950
$ws.Range("J16").Value = 2766.6667
$ws.Range("K16").Value = 10101950
$ws.Range("L16").Value = 2766.6667
$ws.Range("M16").Value = -10101663
$ws.Range("N16").Value = -3340.6667
$ws.Range("H31").Value = 3576.9355
$ws.Range("I31").Value = 1716.1538
$ws.Range("J31").Value = 4920.8335
$ws.Range("K31").Value = 1716.1538
$ws.Range("L31").Value = 4920.8335
$ws.Range("M31").Value = -1421.1538
$ws.Range("N31").Value = -5510.8335
$ws.Range("H34").Value = 3576.9355
$ws.Range("I34").Value = 1716.1538
$ws.Range("J34").Value = 4920.8335
$ws.Range("K34").Value = 1716.1538
$ws.Range("L34").Value = 4920.8335
$ws.Range("M34").Value = -1514.1538
$ws.Range("N34").Value = -5324.8335
$ws.Range("H52").Value = 63375
$ws.Range("J52").Value = 63375
$ws.Range("L52").Value = 63375
$ws.Range("N52").Value = -63963
$ws.Range("H99").Value = 6724
$ws.Range("I99").Value = 4024.4
$ws.Range("J99").Value = 8411.25
$ws.Range("K99").Value = 4024.4
$ws.Range("L99").Value = 8411.25
$ws.Range("M99").Value = -2526.4
$ws.Range("N99").Value = -11407.25
$ws.Range("H105").Value = 1889.2941
$ws.Range("J105").Value = 2625.4285
$ws.Range("L105").Value = 2625.4285
$ws.Range("N105").Value = -6119.4285
$ws.Range("H113").Value = 4832811
$ws.Range("I113").Value = 10101950
$ws.Range("J113").Value = 2766.6667
$ws.Range("K113").Value = 10101950
$ws.Range("L113").Value = 2766.6667
$ws.Range("M113").Value = -10099780
$ws.Range("N113").Value = -7106.6667
$ws.Range("H126").Value = 6724
$ws.Range("I126").Value = 4024.4
$ws.Range("J126").Value = 8411.25
$ws.Range("K126").Value = 12073.2
$ws.Range("L126").Value = 25233.75
$ws.Range("M126").Value = -9603.200000000001
$ws.Range("N126").Value = -30173.75
$ws.Range("H134").Value = 8207.723
$ws.Range("I134").Value = 8903
$ws.Range("K134").Value = 26709
$ws.Range("M134").Value = -24174
$ws.Range("H137").Value = 50580
$ws.Range("J137").Value = 50580
$ws.Range("L137").Value = 50580
$ws.Range("N137").Value = -60780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 10994.667
$ws.Range("I34").Value = 349.33334
$ws.Range("J34").Value = 14543.111
$ws.Range("K34").Value = 1048.00002
$ws.Range("L34").Value = 43629.333
$ws.Range("M34").Value = -964.0000199999999
$ws.Range("N34").Value = -43797.333
$ws.Range("H69").Value = 3310.2222
$ws.Range("I69").Value = 1703
$ws.Range("J69").Value = 4596
$ws.Range("K69").Value = 5109
$ws.Range("L69").Value = 13788
$ws.Range("M69").Value = -4298
$ws.Range("N69").Value = -15410
$ws.Range("H72").Value = 3310.2222
$ws.Range("I72").Value = 1703
$ws.Range("J72").Value = 4596
$ws.Range("K72").Value = 15327
$ws.Range("L72").Value = 41364
$ws.Range("M72").Value = -11271
$ws.Range("N72").Value = -49476
$ws.Range("H113").Value = 902.14813
$ws.Range("I113").Value = 735.36365
$ws.Range("J113").Value = 1636
$ws.Range("K113").Value = 2206.09095
$ws.Range("L113").Value = 4908
$ws.Range("M113").Value = -36.09094999999979
$ws.Range("N113").Value = -9248
$ws.Range("H127").Value = 833.4286
$ws.Range("J127").Value = 833.4286
$ws.Range("L127").Value = 2500.2858
$ws.Range("N127").Value = -12420.2858

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1493.3334
$ws.Range("I113").Value = 1493.3334
$ws.Range("K113").Value = 1493.3334
$ws.Range("M113").Value = 676.6666
$ws.Range("H132").Value = 3146.1853
$ws.Range("I132").Value = 2190.6667
$ws.Range("J132").Value = 3623.9443
$ws.Range("K132").Value = 6572.000100000001
$ws.Range("L132").Value = 10871.8329
$ws.Range("M132").Value = -4042.000100000001
$ws.Range("N132").Value = -15931.8329
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6004.467
$ws.Range("I40").Value = 5506.16
$ws.Range("K40").Value = 5506.16
$ws.Range("M40").Value = -5370.16
$ws.Range("H82").Value = 4915.731
$ws.Range("I82").Value = 5746.45
$ws.Range("J82").Value = 2146.6667
$ws.Range("K82").Value = 5746.45
$ws.Range("L82").Value = 2146.6667
$ws.Range("M82").Value = -5385.45
$ws.Range("N82").Value = -2868.6667
$ws.Range("H85").Value = 4915.731
$ws.Range("I85").Value = 5746.45
$ws.Range("J85").Value = 2146.6667
$ws.Range("K85").Value = 5746.45
$ws.Range("L85").Value = 2146.6667
$ws.Range("M85").Value = -4498.45
$ws.Range("N85").Value = -4642.6667
$ws.Range("H93").Value = 12346588
$ws.Range("I93").Value = 12346588
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 12346588
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -12345340
$ws.Range("N93").ClearContents()
$ws.Range("H106").Value = 25369.666
$ws.Range("J106").Value = 25369.666
$ws.Range("L106").Value = 25369.666
$ws.Range("N106").Value = -27893.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 20150
$ws.Range("J57").Value = 20150
$ws.Range("L57").Value = 20150
$ws.Range("N57").Value = -21658
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H82").Value = 44076.92
$ws.Range("J82").Value = 44076.92
$ws.Range("L82").Value = 44076.92
$ws.Range("N82").Value = -44842.92
$ws.Range("H85").Value = 44076.92
$ws.Range("J85").Value = 44076.92
$ws.Range("L85").Value = 44076.92
$ws.Range("N85").Value = -46728.92
$ws.Range("H99").Value = 30000
$ws.Range("I99").Value = 30000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 30000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -27005
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 880.05884
$ws.Range("I100").Value = 778.875
$ws.Range("K100").Value = 1557.75
$ws.Range("M100").Value = -1016.75
